$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2994
$ws.Range("J17").Value = 2994
$ws.Range("L17").Value = 8982
$ws.Range("N17").Value = -9318

$ws.Range("H62").Value = 5916.6
$ws.Range("I62").Value = 5916.6
$ws.Range("K62").Value = 5916.6
$ws.Range("M62").Value = -5292.6

$ws.Range("H65").Value = 5916.6
$ws.Range("I65").Value = 5916.6
$ws.Range("K65").Value = 29583
$ws.Range("M65").Value = -26463

$ws.Range("H106").Value = 14746.632
$ws.Range("I106").Value = 14824.1875
$ws.Range("K106").Value = 14824.1875
$ws.Range("M106").Value = -14193.1875

$ws.Range("H112").Value = 2437
$ws.Range("I112").Value = 945
$ws.Range("J112").Value = 2561.3333
$ws.Range("K112").Value = 2835
$ws.Range("L112").Value = 7683.999899999999
$ws.Range("M112").Value = -1727
$ws.Range("N112").Value = -9899.999899999999

$ws.Range("H116").Value = 5618.4287
$ws.Range("I116").Value = 5738.1665
$ws.Range("J116").Value = 4900
$ws.Range("K116").Value = 5738.1665
$ws.Range("L116").Value = 4900
$ws.Range("M116").Value = -2296.1665
$ws.Range("N116").Value = -11784

$ws.Range("H138").Value = 4283.1787
$ws.Range("J138").Value = 3762.7222
$ws.Range("L138").Value = 11288.1666
$ws.Range("N138").Value = -21568.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6609.4243
$ws.Range("I32").Value = 4398.5864
$ws.Range("K32").Value = 4398.5864
$ws.Range("M32").Value = -4111.5864

$ws.Range("H45").Value = 998.5
$ws.Range("I45").Value = 998.5
$ws.Range("K45").Value = 998.5
$ws.Range("M45").Value = -621.5

$ws.Range("H74").Value = 1664.091
$ws.Range("I74").Value = 1279.8334
$ws.Range("J74").Value = 3393.25
$ws.Range("K74").Value = 1279.8334
$ws.Range("L74").Value = 3393.25
$ws.Range("M74").Value = -405.8334
$ws.Range("N74").Value = -5141.25

$ws.Range("H77").Value = 1664.091
$ws.Range("I77").Value = 1279.8334
$ws.Range("J77").Value = 3393.25
$ws.Range("K77").Value = 6399.166999999999
$ws.Range("L77").Value = 16966.25
$ws.Range("M77").Value = -2031.166999999999
$ws.Range("N77").Value = -25702.25

$ws.Range("H102").Value = 4650.222
$ws.Range("I102").Value = 4650.222
$ws.Range("K102").Value = 4650.222
$ws.Range("M102").Value = -3028.222

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws.Range("H132").Value = 2955.439
$ws.Range("I132").Value = 2955.439
$ws.Range("K132").Value = 8866.316999999999
$ws.Range("M132").Value = -6336.316999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 17999.334
$ws.Range("I82").Value = 17999.334
$ws.Range("K82").Value = 17999.334
$ws.Range("M82").Value = -17616.334

$ws.Range("H85").Value = 17999.334
$ws.Range("I85").Value = 17999.334
$ws.Range("K85").Value = 17999.334
$ws.Range("M85").Value = -16673.334

$ws.Range("H94").Value = 1433
$ws.Range("I94").Value = 1433
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1433
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -982
$ws.Range("N94").ClearContents()

$ws.Range("H134").Value = 1984.5
$ws.Range("I134").Value = 1981.3334
$ws.Range("K134").Value = 5944.0002
$ws.Range("M134").Value = -3409.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 2000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H27").Value = 2000
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H31").Value = 2275.64
$ws.Range("I31").Value = 1787.6666
$ws.Range("K31").Value = 1787.6666
$ws.Range("M31").Value = -1492.6666

$ws.Range("H34").Value = 2275.64
$ws.Range("I34").Value = 1787.6666
$ws.Range("K34").Value = 1787.6666
$ws.Range("M34").Value = -1585.6666

$ws.Range("H99").Value = 13963.125
$ws.Range("I99").Value = 6314.6665
$ws.Range("J99").Value = 18552.2
$ws.Range("K99").Value = 6314.6665
$ws.Range("L99").Value = 18552.2
$ws.Range("M99").Value = -4816.6665
$ws.Range("N99").Value = -21548.2

$ws.Range("H126").Value = 13963.125
$ws.Range("I126").Value = 6314.6665
$ws.Range("J126").Value = 18552.2
$ws.Range("K126").Value = 18943.9995
$ws.Range("L126").Value = 55656.60000000001
$ws.Range("M126").Value = -16473.9995
$ws.Range("N126").Value = -60596.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3098984.2
$ws.Range("I4").Value = 935154.4
$ws.Range("J4").Value = 15000049
$ws.Range("K4").Value = 2805463.2
$ws.Range("L4").Value = 45000147
$ws.Range("M4").Value = -2805351.2
$ws.Range("N4").Value = -45000371

$ws.Range("H56").Value = 11654.429
$ws.Range("I56").Value = 11654.429
$ws.Range("K56").Value = 11654.429
$ws.Range("M56").Value = -11124.429

$ws.Range("H93").Value = 14562.5
$ws.Range("I93").Value = 750
$ws.Range("K93").Value = 2250
$ws.Range("M93").Value = -378

$ws.Range("H140").Value = 717403.1
$ws.Range("I140").Value = 717403.1
$ws.Range("K140").Value = 2152209.3
$ws.Range("M140").Value = -2147029.3

$ws.Range("H141").Value = 3590.8333
$ws.Range("I141").Value = 3590.8333
$ws.Range("K141").Value = 10772.4999
$ws.Range("M141").Value = -5592.499899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 22500
$ws.Range("I22").Value = 15000
$ws.Range("J22").Value = 30000
$ws.Range("K22").Value = 15000
$ws.Range("L22").Value = 30000
$ws.Range("M22").Value = -14471
$ws.Range("N22").Value = -31058

$ws.Range("H95").Value = 15172
$ws.Range("I95").Value = 10000
$ws.Range("K95").Value = 10000
$ws.Range("M95").Value = -7254

$ws.Range("H107").Value = 2216.1667
$ws.Range("I107").Value = 659.6
$ws.Range("K107").Value = 659.6
$ws.Range("M107").Value = 1260.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 5000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 5000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 5000
$ws.Range("N25").Value = -5460
$ws.Range("M25").ClearContents()

$ws.Range("H40").Value = 71435480
$ws.Range("I40").Value = 166671300
$ws.Range("K40").Value = 166671300
$ws.Range("M40").Value = -166671164

$ws.Range("H55").Value = 1158.8422
$ws.Range("I55").Value = 679
$ws.Range("J55").Value = 1981.4286
$ws.Range("K55").Value = 679
$ws.Range("L55").Value = 1981.4286
$ws.Range("M55").Value = -506
$ws.Range("N55").Value = -2327.4286

$ws.Range("H82").Value = 1109.3125
$ws.Range("I82").Value = 816.375
$ws.Range("J82").Value = 1402.25
$ws.Range("K82").Value = 816.375
$ws.Range("L82").Value = 1402.25
$ws.Range("M82").Value = -455.375
$ws.Range("N82").Value = -2124.25

$ws.Range("H85").Value = 1109.3125
$ws.Range("I85").Value = 816.375
$ws.Range("J85").Value = 1402.25
$ws.Range("K85").Value = 816.375
$ws.Range("L85").Value = 1402.25
$ws.Range("M85").Value = 431.625
$ws.Range("N85").Value = -3898.25

$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

$ws.Range("H132").Value = 11989
$ws.Range("I132").Value = 4386.8
$ws.Range("K132").Value = 13160.4
$ws.Range("M132").Value = -10630.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7267.04
$ws.Range("I81").Value = 4123.933
$ws.Range("J81").Value = 11981.7
$ws.Range("K81").Value = 8247.866
$ws.Range("L81").Value = 23963.4
$ws.Range("M81").Value = -7186.866
$ws.Range("N81").Value = -26085.4

$ws.Range("H84").Value = 7267.04
$ws.Range("I84").Value = 4123.933
$ws.Range("J84").Value = 11981.7
$ws.Range("K84").Value = 41239.33
$ws.Range("L84").Value = 119817
$ws.Range("M84").Value = -35935.33
$ws.Range("N84").Value = -130425

$ws.Range("H100").Value = 1142.9231
$ws.Range("I100").Value = 993.2222
$ws.Range("K100").Value = 1986.4444
$ws.Range("M100").Value = -1445.4444

$ws.Range("H113").Value = 360.5625
$ws.Range("I113").Value = 383.64285
$ws.Range("K113").Value = 1150.92855
$ws.Range("M113").Value = 1019.07145

$ws.Range("H132").Value = 166668210
$ws.Range("I132").Value = 1821.75
$ws.Range("K132").Value = 5465.25
$ws.Range("M132").Value = -2935.25

Write-Host "edit.ps1 applied successfully"
